# "Generate Report for Archive"
#
# 1. Status text "Ready for handoff" -> "In Translation" on every sheet
#    that surfaces it (Overview!E2:F2, zh-cn!C2, de-de!C2).
# 2. Narrow the "Status" columns (Overview E:F, zh-cn C, de-de C) from
#    ~17.22 chars down to ~13.41 chars.

$wb = $excel.ActiveWorkbook

# --- Overview sheet --------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5

# --- zh-cn sheet -------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "In Translation"
$wsZhCn.Columns.Item(3).ColumnWidth = 12.5

# --- de-de sheet -------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "In Translation"
$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
